# Natmi following Dr Hou advice
# Updates the Lgi4-Adam23 LR-pair sheet: adds a new "ECs" / "M1" cluster
# category and rewrites the 6 original data rows (2-7) plus appends 6 new
# data rows (8-13) covering the full 3x4 sending/target cluster grid.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Lgi4"
$ws.Range("C2").Value = "Adam23"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.06754433333333333
$ws.Range("H2").Value = 0.202633
$ws.Range("I2").Value = 0.02266275549884949
$ws.Range("J2").Value = 0.02266275549884949
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.4998576666666667
$ws.Range("N2").Value = 1.499573
$ws.Range("O2").Value = 0.02650617333988447
$ws.Range("P2").Value = 0.02650617333988446
$ws.Range("Q2").Value = 0.03376255285655556
$ws.Range("R2").Value = 0.303862975709
$ws.Range("S2").Value = 0.0006007029256119244
$ws.Range("T2").Value = 0.0006007029256119243

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Lgi4"
$ws.Range("C3").Value = "Adam23"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.06754433333333333
$ws.Range("H3").Value = 0.202633
$ws.Range("I3").Value = 0.02266275549884949
$ws.Range("J3").Value = 0.02266275549884949
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 17.63262733333333
$ws.Range("N3").Value = 52.897882
$ws.Range("O3").Value = 0.9350131201380354
$ws.Range("P3").Value = 0.9350131201380353
$ws.Range("Q3").Value = 1.190984058145111
$ws.Range("R3").Value = 10.718856523306
$ws.Range("S3").Value = 0.02118997372990468
$ws.Range("T3").Value = 0.02118997372990468

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Lgi4"
$ws.Range("C4").Value = "Adam23"
$ws.Range("D4").Value = "M1"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.06754433333333333
$ws.Range("H4").Value = 0.202633
$ws.Range("I4").Value = 0.02266275549884949
$ws.Range("J4").Value = 0.02266275549884949
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.01414533333333333
$ws.Range("N4").Value = 0.042436
$ws.Range("O4").Value = 0.0007500908404267997
$ws.Range("P4").Value = 0.0007500908404267996
$ws.Range("Q4").Value = 0.0009554371097777779
$ws.Range("R4").Value = 0.008598933988
$ws.Range("S4").Value = 0.00001699912531851909
$ws.Range("T4").Value = 0.00001699912531851909

# Row 5
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Lgi4"
$ws.Range("C5").Value = "Adam23"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.06754433333333333
$ws.Range("H5").Value = 0.202633
$ws.Range("I5").Value = 0.02266275549884949
$ws.Range("J5").Value = 0.02266275549884949
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.7115299999999999
$ws.Range("N5").Value = 2.13459
$ws.Range("O5").Value = 0.03773061568165336
$ws.Range("P5").Value = 0.03773061568165335
$ws.Range("Q5").Value = 0.04805981949666666
$ws.Range("R5").Value = 0.43253837547
$ws.Range("S5").Value = 0.0008550797180143663
$ws.Range("T5").Value = 0.0008550797180143663

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Lgi4"
$ws.Range("C6").Value = "Adam23"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 2.302381666666667
$ws.Range("H6").Value = 6.907145
$ws.Range("I6").Value = 0.7725046677002302
$ws.Range("J6").Value = 0.7725046677002302
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.4998576666666667
$ws.Range("N6").Value = 1.499573
$ws.Range("O6").Value = 0.02650617333988447
$ws.Range("P6").Value = 0.02650617333988446
$ws.Range("Q6").Value = 1.150863127676111
$ws.Range("R6").Value = 10.357768149085
$ws.Range("S6").Value = 0.02047614262793215
$ws.Range("T6").Value = 0.02047614262793215

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Lgi4"
$ws.Range("C7").Value = "Adam23"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 2.302381666666667
$ws.Range("H7").Value = 6.907145
$ws.Range("I7").Value = 0.7725046677002302
$ws.Range("J7").Value = 0.7725046677002302
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 17.63262733333333
$ws.Range("N7").Value = 52.897882
$ws.Range("O7").Value = 0.9350131201380354
$ws.Range("P7").Value = 0.9350131201380353
$ws.Range("Q7").Value = 40.59703790743222
$ws.Range("R7").Value = 365.37334116689
$ws.Range("S7").Value = 0.7223019996675885
$ws.Range("T7").Value = 0.7223019996675883

# Row 8
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Lgi4"
$ws.Range("C8").Value = "Adam23"
$ws.Range("D8").Value = "M1"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 2.302381666666667
$ws.Range("H8").Value = 6.907145
$ws.Range("I8").Value = 0.7725046677002302
$ws.Range("J8").Value = 0.7725046677002302
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.01414533333333333
$ws.Range("N8").Value = 0.042436
$ws.Range("O8").Value = 0.0007500908404267997
$ws.Range("P8").Value = 0.0007500908404267996
$ws.Range("Q8").Value = 0.03256795613555556
$ws.Range("R8").Value = 0.29311160522
$ws.Range("S8").Value = 0.0005794486754288913
$ws.Range("T8").Value = 0.0005794486754288912

# Row 9
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Lgi4"
$ws.Range("C9").Value = "Adam23"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 2.302381666666667
$ws.Range("H9").Value = 6.907145
$ws.Range("I9").Value = 0.7725046677002302
$ws.Range("J9").Value = 0.7725046677002302
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.7115299999999999
$ws.Range("N9").Value = 2.13459
$ws.Range("O9").Value = 0.03773061568165336
$ws.Range("P9").Value = 0.03773061568165335
$ws.Range("Q9").Value = 1.638213627283333
$ws.Range("R9").Value = 14.74392264555
$ws.Range("S9").Value = 0.02914707672928072
$ws.Range("T9").Value = 0.02914707672928072

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Lgi4"
$ws.Range("C10").Value = "Adam23"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.6104853333333333
$ws.Range("H10").Value = 1.831456
$ws.Range("I10").Value = 0.2048325768009203
$ws.Range("J10").Value = 0.2048325768009203
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.4998576666666667
$ws.Range("N10").Value = 1.499573
$ws.Range("O10").Value = 0.02650617333988447
$ws.Range("P10").Value = 0.02650617333988446
$ws.Range("Q10").Value = 0.3051557742542222
$ws.Range("R10").Value = 2.746401968288
$ws.Range("S10").Value = 0.005429327786340392
$ws.Range("T10").Value = 0.005429327786340391

# Row 11
$ws.Range("A11").Value = "sCs"
$ws.Range("B11").Value = "Lgi4"
$ws.Range("C11").Value = "Adam23"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.6104853333333333
$ws.Range("H11").Value = 1.831456
$ws.Range("I11").Value = 0.2048325768009203
$ws.Range("J11").Value = 0.2048325768009203
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 17.63262733333333
$ws.Range("N11").Value = 52.897882
$ws.Range("O11").Value = 0.9350131201380354
$ws.Range("P11").Value = 0.9350131201380353
$ws.Range("Q11").Value = 10.76446037513244
$ws.Range("R11").Value = 96.880143376192
$ws.Range("S11").Value = 0.1915211467405422
$ws.Range("T11").Value = 0.1915211467405422

# Row 12
$ws.Range("A12").Value = "sCs"
$ws.Range("B12").Value = "Lgi4"
$ws.Range("C12").Value = "Adam23"
$ws.Range("D12").Value = "M1"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 0.6104853333333333
$ws.Range("H12").Value = 1.831456
$ws.Range("I12").Value = 0.2048325768009203
$ws.Range("J12").Value = 0.2048325768009203
$ws.Range("K12").Value = 1
$ws.Range("L12").Value = 0.3333333333333333
$ws.Range("M12").Value = 0.01414533333333333
$ws.Range("N12").Value = 0.042436
$ws.Range("O12").Value = 0.0007500908404267997
$ws.Range("P12").Value = 0.0007500908404267996
$ws.Range("Q12").Value = 0.008635518535111112
$ws.Range("R12").Value = 0.077719666816
$ws.Range("S12").Value = 0.0001536430396793893
$ws.Range("T12").Value = 0.0001536430396793893

# Row 13
$ws.Range("A13").Value = "sCs"
$ws.Range("B13").Value = "Lgi4"
$ws.Range("C13").Value = "Adam23"
$ws.Range("D13").Value = "sCs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 0.6104853333333333
$ws.Range("H13").Value = 1.831456
$ws.Range("I13").Value = 0.2048325768009203
$ws.Range("J13").Value = 0.2048325768009203
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.7115299999999999
$ws.Range("N13").Value = 2.13459
$ws.Range("O13").Value = 0.03773061568165336
$ws.Range("P13").Value = 0.03773061568165335
$ws.Range("Q13").Value = 0.4343786292266666
$ws.Range("R13").Value = 3.90940766304
$ws.Range("S13").Value = 0.00772845923435827
$ws.Range("T13").Value = 0.007728459234358269
